$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for the Col2a1-Mag ligand-receptor pair table.
# Each assignment below corresponds to a recomputed metric cell from the
# refreshed TPM expression data (commit: "update scripts wuth new tpm").
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1581976666666667
$ws.Range("H2").Value = 0.474593
$ws.Range("I2").Value = 0.1400666049254827
$ws.Range("J2").Value = 0.1400666049254826
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8180823333333334
$ws.Range("N2").Value = 2.454247
$ws.Range("O2").Value = 0.5115352725808422
$ws.Range("P2").Value = 0.5115352725808422
$ws.Range("Q2").Value = 0.1294187162745556
$ws.Range("R2").Value = 1.164768446471
$ws.Range("S2").Value = 0.07164900893002991
$ws.Range("T2").Value = 0.0716490089300299
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1581976666666667
$ws.Range("H3").Value = 0.474593
$ws.Range("I3").Value = 0.1400666049254827
$ws.Range("J3").Value = 0.1400666049254826
$ws.Range("M3").Value = 0.7811863333333333
$ws.Range("N3").Value = 2.343559
$ws.Range("O3").Value = 0.4884647274191578
$ws.Range("P3").Value = 0.4884647274191579
$ws.Range("Q3").Value = 0.1235818551652222
$ws.Range("R3").Value = 1.112236696487
$ws.Range("S3").Value = 0.06841759599545276
$ws.Range("T3").Value = 0.06841759599545276
$ws.Range("I4").Value = 0.7029419733214338
$ws.Range("J4").Value = 0.7029419733214337
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8180823333333334
$ws.Range("N4").Value = 2.454247
$ws.Range("O4").Value = 0.5115352725808422
$ws.Range("P4").Value = 0.5115352725808422
$ws.Range("Q4").Value = 0.649504197315
$ws.Range("R4").Value = 5.845537775835
$ws.Range("S4").Value = 0.3595796139314947
$ws.Range("T4").Value = 0.3595796139314947
$ws.Range("I5").Value = 0.7029419733214338
$ws.Range("J5").Value = 0.7029419733214337
$ws.Range("M5").Value = 0.7811863333333333
$ws.Range("N5").Value = 2.343559
$ws.Range("O5").Value = 0.4884647274191578
$ws.Range("P5").Value = 0.4884647274191579
$ws.Range("Q5").Value = 0.6202111715549999
$ws.Range("R5").Value = 5.581900543994999
$ws.Range("S5").Value = 0.3433623593899391
$ws.Range("T5").Value = 0.343362359389939
$ws.Range("G6").Value = 0.1199896666666667
$ws.Range("H6").Value = 0.359969
$ws.Range("I6").Value = 0.106237630366274
$ws.Range("J6").Value = 0.106237630366274
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8180823333333334
$ws.Range("N6").Value = 2.454247
$ws.Range("O6").Value = 0.5115352725808422
$ws.Range("P6").Value = 0.5115352725808422
$ws.Range("Q6").Value = 0.09816142648255555
$ws.Range("R6").Value = 0.883452838343
$ws.Range("S6").Value = 0.05434429520775472
$ws.Range("T6").Value = 0.0543442952077547
$ws.Range("G7").Value = 0.1199896666666667
$ws.Range("H7").Value = 0.359969
$ws.Range("I7").Value = 0.106237630366274
$ws.Range("J7").Value = 0.106237630366274
$ws.Range("M7").Value = 0.7811863333333333
$ws.Range("N7").Value = 2.343559
$ws.Range("O7").Value = 0.4884647274191578
$ws.Range("P7").Value = 0.4884647274191579
$ws.Range("Q7").Value = 0.09373428774122222
$ws.Range("R7").Value = 0.8436085896709999
$ws.Range("S7").Value = 0.05189333515851927
$ws.Range("T7").Value = 0.05189333515851926
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.05732366666666666
$ws.Range("H8").Value = 0.171971
$ws.Range("I8").Value = 0.05075379138680971
$ws.Range("J8").Value = 0.05075379138680969
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8180823333333334
$ws.Range("N8").Value = 2.454247
$ws.Range("O8").Value = 0.5115352725808422
$ws.Range("P8").Value = 0.5115352725808422
$ws.Range("Q8").Value = 0.04689547898188889
$ws.Range("R8").Value = 0.422059310837
$ws.Range("S8").Value = 0.02596235451156291
$ws.Range("T8").Value = 0.0259623545115629
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.05732366666666666
$ws.Range("H9").Value = 0.171971
$ws.Range("I9").Value = 0.05075379138680971
$ws.Range("J9").Value = 0.05075379138680969
$ws.Range("M9").Value = 0.7811863333333333
$ws.Range("N9").Value = 2.343559
$ws.Range("O9").Value = 0.4884647274191578
$ws.Range("P9").Value = 0.4884647274191579
$ws.Range("Q9").Value = 0.04478046497655555
$ws.Range("R9").Value = 0.4030241847889999
$ws.Range("S9").Value = 0.0247914368752468
$ws.Range("T9").Value = 0.0247914368752468
